$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Add new Quina draws (rows 342-374) below the existing data table.
# ---------------------------------------------------------------

# Step 1: copy the currently-last styled block (A337:F341) and use Insert
# (shift cells down) to stamp out 21 new rows that carry the same cell
# style (these become rows 354-374 once the 12 rows below are inserted in
# step 3).
$ws.Range("A337:F341").Copy()
$ws.Range("A342:F346").Insert(-4121)
$ws.Range("A342:F346").Insert(-4121)
$ws.Range("A342:F346").Insert(-4121)
$ws.Range("A342:F346").Insert(-4121)
$ws.Range("A342:F342").Insert(-4121)
$excel.CutCopyMode = 0

# Step 2: the old "last page" rows (337-341) no longer sit at the bottom of
# the table, so drop their now-redundant highlight style, matching the look
# of the other regular data rows.
$ws.Range("A337:F341").ClearFormats()

# Step 3: insert 12 unstyled rows right after row 341 (they inherit the "no
# style" formatting of row 341 above them) to hold draws 6833-6844.
$ws.Range("A342:F353").Insert(-4121)

# Step 4: populate the 12 new unstyled rows (342-353).
$ws.Range("A342").Value = 6833
$ws.Range("B342").Value = 3
$ws.Range("C342").Value = 7
$ws.Range("D342").Value = 9
$ws.Range("E342").Value = 22
$ws.Range("F342").Value = 73

$ws.Range("A343").Value = 6834
$ws.Range("B343").Value = 11
$ws.Range("C343").Value = 41
$ws.Range("D343").Value = 51
$ws.Range("E343").Value = 78
$ws.Range("F343").Value = 79

$ws.Range("A344").Value = 6835
$ws.Range("B344").Value = 4
$ws.Range("C344").Value = 13
$ws.Range("D344").Value = 30
$ws.Range("E344").Value = 49
$ws.Range("F344").Value = 68

$ws.Range("A345").Value = 6836
$ws.Range("B345").Value = 5
$ws.Range("C345").Value = 16
$ws.Range("D345").Value = 29
$ws.Range("E345").Value = 33
$ws.Range("F345").Value = 76

$ws.Range("A346").Value = 6837
$ws.Range("B346").Value = 13
$ws.Range("C346").Value = 25
$ws.Range("D346").Value = 57
$ws.Range("E346").Value = 64
$ws.Range("F346").Value = 75

$ws.Range("A347").Value = 6838
$ws.Range("B347").Value = 31
$ws.Range("C347").Value = 33
$ws.Range("D347").Value = 50
$ws.Range("E347").Value = 70
$ws.Range("F347").Value = 77

$ws.Range("A348").Value = 6839
$ws.Range("B348").Value = 11
$ws.Range("C348").Value = 26
$ws.Range("D348").Value = 28
$ws.Range("E348").Value = 44
$ws.Range("F348").Value = 61

$ws.Range("A349").Value = 6840
$ws.Range("B349").Value = 26
$ws.Range("C349").Value = 36
$ws.Range("D349").Value = 43
$ws.Range("E349").Value = 46
$ws.Range("F349").Value = 74

$ws.Range("A350").Value = 6841
$ws.Range("B350").Value = 12
$ws.Range("C350").Value = 25
$ws.Range("D350").Value = 33
$ws.Range("E350").Value = 42
$ws.Range("F350").Value = 74

$ws.Range("A351").Value = 6842
$ws.Range("B351").Value = 12
$ws.Range("C351").Value = 25
$ws.Range("D351").Value = 33
$ws.Range("E351").Value = 41
$ws.Range("F351").Value = 52

$ws.Range("A352").Value = 6843
$ws.Range("B352").Value = 4
$ws.Range("C352").Value = 14
$ws.Range("D352").Value = 45
$ws.Range("E352").Value = 71
$ws.Range("F352").Value = 80

$ws.Range("A353").Value = 6844
$ws.Range("B353").Value = 5
$ws.Range("C353").Value = 9
$ws.Range("D353").Value = 11
$ws.Range("E353").Value = 52
$ws.Range("F353").Value = 59

# Step 5: populate the 21 new styled rows (354-374), overwriting the copied
# placeholder values from step 1 with the real draw numbers.
$ws.Range("A354").Value = 6845
$ws.Range("B354").Value = 30
$ws.Range("C354").Value = 45
$ws.Range("D354").Value = 56
$ws.Range("E354").Value = 57
$ws.Range("F354").Value = 62

$ws.Range("A355").Value = 6846
$ws.Range("B355").Value = 4
$ws.Range("C355").Value = 49
$ws.Range("D355").Value = 51
$ws.Range("E355").Value = 57
$ws.Range("F355").Value = 63

$ws.Range("A356").Value = 6847
$ws.Range("B356").Value = 17
$ws.Range("C356").Value = 18
$ws.Range("D356").Value = 27
$ws.Range("E356").Value = 66
$ws.Range("F356").Value = 71

$ws.Range("A357").Value = 6848
$ws.Range("B357").Value = 31
$ws.Range("C357").Value = 32
$ws.Range("D357").Value = 34
$ws.Range("E357").Value = 45
$ws.Range("F357").Value = 80

$ws.Range("A358").Value = 6849
$ws.Range("B358").Value = 13
$ws.Range("C358").Value = 22
$ws.Range("D358").Value = 23
$ws.Range("E358").Value = 30
$ws.Range("F358").Value = 67

$ws.Range("A359").Value = 6850
$ws.Range("B359").Value = 34
$ws.Range("C359").Value = 55
$ws.Range("D359").Value = 61
$ws.Range("E359").Value = 71
$ws.Range("F359").Value = 72

$ws.Range("A360").Value = 6851
$ws.Range("B360").Value = 13
$ws.Range("C360").Value = 18
$ws.Range("D360").Value = 32
$ws.Range("E360").Value = 42
$ws.Range("F360").Value = 55

$ws.Range("A361").Value = 6852
$ws.Range("B361").Value = 24
$ws.Range("C361").Value = 46
$ws.Range("D361").Value = 50
$ws.Range("E361").Value = 62
$ws.Range("F361").Value = 68

$ws.Range("A362").Value = 6853
$ws.Range("B362").Value = 7
$ws.Range("C362").Value = 9
$ws.Range("D362").Value = 48
$ws.Range("E362").Value = 54
$ws.Range("F362").Value = 75

$ws.Range("A363").Value = 6854
$ws.Range("B363").Value = 8
$ws.Range("C363").Value = 42
$ws.Range("D363").Value = 53
$ws.Range("E363").Value = 65
$ws.Range("F363").Value = 68

$ws.Range("A364").Value = 6855
$ws.Range("B364").Value = 4
$ws.Range("C364").Value = 6
$ws.Range("D364").Value = 33
$ws.Range("E364").Value = 56
$ws.Range("F364").Value = 63

$ws.Range("A365").Value = 6856
$ws.Range("B365").Value = 5
$ws.Range("C365").Value = 10
$ws.Range("D365").Value = 16
$ws.Range("E365").Value = 33
$ws.Range("F365").Value = 53

$ws.Range("A366").Value = 6857
$ws.Range("B366").Value = 21
$ws.Range("C366").Value = 32
$ws.Range("D366").Value = 34
$ws.Range("E366").Value = 41
$ws.Range("F366").Value = 46

$ws.Range("A367").Value = 6858
$ws.Range("B367").Value = 31
$ws.Range("C367").Value = 55
$ws.Range("D367").Value = 70
$ws.Range("E367").Value = 78
$ws.Range("F367").Value = 80

$ws.Range("A368").Value = 6859
$ws.Range("B368").Value = 23
$ws.Range("C368").Value = 33
$ws.Range("D368").Value = 36
$ws.Range("E368").Value = 39
$ws.Range("F368").Value = 79

$ws.Range("A369").Value = 6860
$ws.Range("B369").Value = 12
$ws.Range("C369").Value = 46
$ws.Range("D369").Value = 58
$ws.Range("E369").Value = 62
$ws.Range("F369").Value = 69

$ws.Range("A370").Value = 6861
$ws.Range("B370").Value = 26
$ws.Range("C370").Value = 32
$ws.Range("D370").Value = 42
$ws.Range("E370").Value = 57
$ws.Range("F370").Value = 78

$ws.Range("A371").Value = 6862
$ws.Range("B371").Value = 11
$ws.Range("C371").Value = 39
$ws.Range("D371").Value = 51
$ws.Range("E371").Value = 62
$ws.Range("F371").Value = 68

$ws.Range("A372").Value = 6863
$ws.Range("B372").Value = 28
$ws.Range("C372").Value = 47
$ws.Range("D372").Value = 50
$ws.Range("E372").Value = 65
$ws.Range("F372").Value = 79

$ws.Range("A373").Value = 6864
$ws.Range("B373").Value = 6
$ws.Range("C373").Value = 18
$ws.Range("D373").Value = 55
$ws.Range("E373").Value = 58
$ws.Range("F373").Value = 78

$ws.Range("A374").Value = 6865
$ws.Range("B374").Value = 29
$ws.Range("C374").Value = 42
$ws.Range("D374").Value = 55
$ws.Range("E374").Value = 71
$ws.Range("F374").Value = 75

# Step 6: scroll the view down and select the newly added last block, as the
# author did after typing in the new draws.
$ws.Range("B354:F374").Select()
$excel.ActiveWindow.ScrollRow = 332
$excel.ActiveWindow.ScrollColumn = 1
